$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("A2").Value = 170
$ws.Range("B2").Value = 40000
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "nothing"

# Row 3 updates
$ws.Range("A3").Value = 180
$ws.Range("B3").Value = 50000
$ws.Range("H3").Value = "nothing"
$ws.Range("K3").Value = 1

# Update selection to H7
$ws.Range("H7").Select()
